# Add a new "Row bootstrap class css" issue row (B3/C3) and a new
# "Countries, Cities, CityAreas Tables" issue row at the bottom (A28),
# plus small formatting touch-ups (center alignment on column B y/n cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix capitalisation of the existing "row bootstrap class css" item ---
$ws.Range("A17").Value = "Row bootstrap class css"

# --- New issue row: "Row bootstrap class css" / Resolved=y / Remarks ---
$ws.Range("B3").Value = "y"
$ws.Range("C3").Value = "We need to speicify order of the element to  take it to first like order-first class in bootstrap"

# --- Center-align the "Resolved" (y) cells in column B ---
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("B5").HorizontalAlignment = -4108

# --- Widen column C to fit the new, much longer Remarks text ---
$ws.Columns.Item(3).ColumnWidth = 84.140625

# --- New issue row at the bottom ---
$ws.Range("A28").Value = "Countries, Cities, CityAreas Tables"

# --- Selection / view matches the author's final state ---
$ws.Range("C5").Select()
